# Applies the weekly cryptos price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.271.41'
$ws.Range('E2').Value = '  -2.94%  '
$ws.Range('D3').Value = '2.997.84'
$ws.Range('E3').Value = '  -3.64%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''580.64'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.91%  '
$ws.Range('D6').Value = '''146.47'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -6.57%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -3.32%  '
$ws.Range('D9').Value = '2.999.68'
$ws.Range('E9').Value = '  -3.71%  '
$ws.Range('D10').Value = '''0.148'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -6.36%  '
$ws.Range('D11').Value = '''5.63'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -4.82%  '
$ws.Range('E12').Value = '  -2.62%  '
$ws.Range('E13').Value = '  -5.05%  '
$ws.Range('D14').Value = '''34.53'
$ws.Range('D14').Style = "Normal"
$ws.Range('E15').Value = '  +1.78%  '
$ws.Range('D16').Value = '3.499.09'
$ws.Range('E16').Value = '  -3.50%  '
$ws.Range('D17').Value = '''7.09'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.95%  '
$ws.Range('D18').Value = '62.323.71'
$ws.Range('E18').Value = '  -2.78%  '
$ws.Range('D19').Value = '3.003.99'
$ws.Range('E19').Value = '  -3.51%  '
$ws.Range('D20').Value = '''453.39'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -5.46%  '
$ws.Range('D21').Value = '''13.85'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -4.42%  '
$ws.Range('D22').Value = '''0.678'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -4.77%  '
$ws.Range('D23').Value = '''7.29'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -3.97%  '
$ws.Range('D24').Value = '''80.07'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.47%  '
$ws.Range('D25').Value = '''2.27'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -6.89%  '
$ws.Range('D26').Value = '''12.27'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -5.24%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').Value = '''0.999'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.11%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').Value = '''9.99'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -3.79%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '''2.61'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.91%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').Value = '''7.13'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -4.46%  '
$ws.Range('E32').Value = '  -4.95%  '
$ws.Range('D33').Value = '''26.86'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -2.11%  '
$ws.Range('E34').Value = '  -5.62%  '
$ws.Range('D35').Value = '''1.02'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -3.57%  '
$ws.Range('D36').Value = '0.0₃0790'
$ws.Range('E36').Value = '  -5.94%  '
$ws.Range('D37').Value = '''5.73'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -4.96%  '
$ws.Range('D38').Value = '''2.11'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -5.83%  '
$ws.Range('D39').Value = '''50.07'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -2.04%  '
$ws.Range('D40').Value = '''9.03'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.03%  '
$ws.Range('E41').Value = '  -10.89%  '
$ws.Range('D42').Value = '''409.33'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -7.04%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').Value = '''0.111'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.18%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').Value = '''0.275'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -5.54%  '
$ws.Range('D45').Value = '2.771.17'
$ws.Range('E45').Value = '  -2.42%  '
$ws.Range('D46').Value = '''0.0351'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -3.75%  '
$ws.Range('D47').Value = '''38.03'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -5.39%  '
$ws.Range('D48').Value = '''127.99'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.58%  '
$ws.Range('E50').Value = '  -2.19%  '
$ws.Range('D51').Value = '''23.77'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -6.33%  '
